$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lessons")

$ws.Range("A52").Value = "Hétfő"
$ws.Range("B52").Value = 1
$ws.Range("C52").Value = "13.B"
$ws.Range("D52").Value = "Backend"
$ws.Range("E52").Value = "HP"
$ws.Range("F52").Value = 22
$ws.Range("G52").Value = "mat2"

$ws.Range("A53").Value = "Hétfő"
$ws.Range("B53").Value = 2
$ws.Range("C53").Value = "13.B"
$ws.Range("D53").Value = "Frontend"
$ws.Range("E53").Value = "ÁA"
$ws.Range("F53").Value = 23
$ws.Range("G53").Value = "mat1"

$ws.Range("A54").Value = "Hétfő"
$ws.Range("B54").Value = 2
$ws.Range("C54").Value = "13.B"
$ws.Range("D54").Value = "Backend"
$ws.Range("E54").Value = "HP"
$ws.Range("F54").Value = 22
$ws.Range("G54").Value = "mat2"

$ws.Range("A55").Value = "Hétfő"
$ws.Range("B55").Value = 3
$ws.Range("C55").Value = "13.B"
$ws.Range("D55").Value = "Frontend"
$ws.Range("E55").Value = "ÁA"
$ws.Range("F55").Value = 23
$ws.Range("G55").Value = "mat1"

$ws.Range("A56").Value = "Hétfő"
$ws.Range("B56").Value = 3
$ws.Range("C56").Value = "13.B"
$ws.Range("D56").Value = "Backend"
$ws.Range("E56").Value = "HP"
$ws.Range("F56").Value = 22
$ws.Range("G56").Value = "mat2"

$ws.Range("A57").Value = "Hétfő"
$ws.Range("B57").Value = 4
$ws.Range("C57").Value = "13.B"
$ws.Range("D57").Value = "Frontend"
$ws.Range("E57").Value = "ÁA"
$ws.Range("F57").Value = 23
$ws.Range("G57").Value = "mat1"

$ws.Range("A58").Value = "Hétfő"
$ws.Range("B58").Value = 4
$ws.Range("C58").Value = "13.B"
$ws.Range("D58").Value = "Backend"
$ws.Range("E58").Value = "HP"
$ws.Range("F58").Value = 22
$ws.Range("G58").Value = "mat2"

$ws.Range("A59").Value = "Hétfő"
$ws.Range("B59").Value = 5
$ws.Range("C59").Value = "13.B"
$ws.Range("D59").Value = "Backend"
$ws.Range("E59").Value = "HZs"
$ws.Range("F59").Value = 23
$ws.Range("G59").Value = "mat1"

$ws.Range("A60").Value = "Hétfő"
$ws.Range("B60").Value = 5
$ws.Range("C60").Value = "13.B"
$ws.Range("D60").Value = "Aalkf."
$ws.Range("E60").Value = "VT"
$ws.Range("F60").Value = 30
$ws.Range("G60").Value = "mat2"

$ws.Range("A61").Value = "Hétfő"
$ws.Range("B61").Value = 6
$ws.Range("C61").Value = "13.B"
$ws.Range("D61").Value = "Backend"
$ws.Range("E61").Value = "HZs"
$ws.Range("F61").Value = 23
$ws.Range("G61").Value = "mat1"

$ws.Range("A62").Value = "Hétfő"
$ws.Range("B62").Value = 6
$ws.Range("C62").Value = "13.B"
$ws.Range("D62").Value = "Aalkf."
$ws.Range("E62").Value = "VT"
$ws.Range("F62").Value = 30
$ws.Range("G62").Value = "mat2"

$ws.Range("A63").Value = "Hétfő"
$ws.Range("B63").Value = 7
$ws.Range("C63").Value = "13.B"
$ws.Range("D63").Value = "Backend"
$ws.Range("E63").Value = "HZs"
$ws.Range("F63").Value = 23
$ws.Range("G63").Value = "mat1"

$ws.Range("A64").Value = "Hétfő"
$ws.Range("B64").Value = 7
$ws.Range("C64").Value = "13.B"
$ws.Range("D64").Value = "Aalkf."
$ws.Range("E64").Value = "VT"
$ws.Range("F64").Value = 30
$ws.Range("G64").Value = "mat2"

$ws.Range("A65").Value = "Kedd"
$ws.Range("B65").Value = 1
$ws.Range("C65").Value = "13.B"
$ws.Range("D65").Value = "Backend"
$ws.Range("E65").Value = "HZs"
$ws.Range("F65").Value = 23
$ws.Range("G65").Value = "mat1"

$ws.Range("A66").Value = "Kedd"
$ws.Range("B66").Value = 1
$ws.Range("C66").Value = "13.B"
$ws.Range("D66").Value = "Frontend"
$ws.Range("E66").Value = "HP"
$ws.Range("F66").Value = 22
$ws.Range("G66").Value = "mat2"

$ws.Range("A67").Value = "Kedd"
$ws.Range("B67").Value = 2
$ws.Range("C67").Value = "13.B"
$ws.Range("D67").Value = "Backend"
$ws.Range("E67").Value = "HZs"
$ws.Range("F67").Value = 23
$ws.Range("G67").Value = "mat1"

$ws.Range("A68").Value = "Kedd"
$ws.Range("B68").Value = 2
$ws.Range("C68").Value = "13.B"
$ws.Range("D68").Value = "Frontend"
$ws.Range("E68").Value = "HP"
$ws.Range("F68").Value = 22
$ws.Range("G68").Value = "mat2"

$ws.Range("A69").Value = "Kedd"
$ws.Range("B69").Value = 3
$ws.Range("C69").Value = "13.B"
$ws.Range("D69").Value = "Backend"
$ws.Range("E69").Value = "HZs"
$ws.Range("F69").Value = 23
$ws.Range("G69").Value = "mat1"

$ws.Range("A70").Value = "Kedd"
$ws.Range("B70").Value = 3
$ws.Range("C70").Value = "13.B"
$ws.Range("D70").Value = "Frontend"
$ws.Range("E70").Value = "HP"
$ws.Range("F70").Value = 22
$ws.Range("G70").Value = "mat2"

$ws.Range("A71").Value = "Kedd"
$ws.Range("B71").Value = 4
$ws.Range("C71").Value = "13.B"
$ws.Range("D71").Value = "Mii nyelv"
$ws.Range("E71").Value = "Iv"
$ws.Range("F71").Value = 21
$ws.Range("G71").Value = "Teljes osztály"

$ws.Range("A72").Value = "Kedd"
$ws.Range("B72").Value = 5
$ws.Range("C72").Value = "13.B"
$ws.Range("D72").Value = "Érettségi felkészítő"
$ws.Range("E72").Value = "TDZs"
$ws.Range("F72").Value = 21
$ws.Range("G72").Value = "mat1"

$ws.Range("A73").Value = "Kedd"
$ws.Range("B73").Value = 6
$ws.Range("C73").Value = "13.B"
$ws.Range("D73").Value = "Érettségi felkészítő"
$ws.Range("E73").Value = "TDZs"
$ws.Range("F73").Value = 21
$ws.Range("G73").Value = "mat1"

$ws.Range("A74").Value = "Kedd"
$ws.Range("B74").Value = 7
$ws.Range("C74").Value = "13.B"
$ws.Range("D74").Value = "Érettségi felkészítő"
$ws.Range("E74").Value = "TDZs"
$ws.Range("F74").Value = 21
$ws.Range("G74").Value = "mat1"

$ws.Range("A75").Value = "Szerda"
$ws.Range("B75").Value = 2
$ws.Range("C75").Value = "13.B"
$ws.Range("D75").Value = "Hittan"
$ws.Range("E75").Value = "BP"
$ws.Range("F75").Value = 9
$ws.Range("G75").Value = "Teljes osztály"

$ws.Range("A76").Value = "Szerda"
$ws.Range("B76").Value = 3
$ws.Range("C76").Value = "13.B"
$ws.Range("D76").Value = "Köznevelés"
$ws.Range("E76").Value = "BSz"
$ws.Range("F76").Value = 20
$ws.Range("G76").Value = "Teljes osztály"

$ws.Range("A77").Value = "Szerda"
$ws.Range("B77").Value = 4
$ws.Range("C77").Value = "13.B"
$ws.Range("D77").Value = "Aalkf."
$ws.Range("E77").Value = "VT"
$ws.Range("F77").Value = 30
$ws.Range("G77").Value = "mat1"

$ws.Range("A78").Value = "Szerda"
$ws.Range("B78").Value = 4
$ws.Range("C78").Value = "13.B"
$ws.Range("D78").Value = "Backend"
$ws.Range("E78").Value = "HP"
$ws.Range("F78").Value = 22
$ws.Range("G78").Value = "mat2"

$ws.Range("A79").Value = "Szerda"
$ws.Range("B79").Value = 5
$ws.Range("C79").Value = "13.B"
$ws.Range("D79").Value = "Aalkf."
$ws.Range("E79").Value = "VT"
$ws.Range("F79").Value = 30
$ws.Range("G79").Value = "mat1"

$ws.Range("A80").Value = "Szerda"
$ws.Range("B80").Value = 5
$ws.Range("C80").Value = "13.B"
$ws.Range("D80").Value = "Backend"
$ws.Range("E80").Value = "HP"
$ws.Range("F80").Value = 22
$ws.Range("G80").Value = "mat2"

$ws.Range("A81").Value = "Szerda"
$ws.Range("B81").Value = 6
$ws.Range("C81").Value = "13.B"
$ws.Range("D81").Value = "Mii nyelv"
$ws.Range("E81").Value = "KKr"
$ws.Range("F81").Value = 17
$ws.Range("G81").Value = "Teljes osztály"

$ws.Range("A82").Value = "Szerda"
$ws.Range("B82").Value = 7
$ws.Range("C82").Value = "13.B"
$ws.Range("D82").Value = "Adatbázis-kezelés"
$ws.Range("E82").Value = "HZ"
$ws.Range("F82").Value = 23
$ws.Range("G82").Value = "mat1"

$ws.Range("A83").Value = "Szerda"
$ws.Range("B83").Value = 7
$ws.Range("C83").Value = "13.B"
$ws.Range("D83").Value = "Aalkf."
$ws.Range("E83").Value = "VT"
$ws.Range("F83").Value = 30
$ws.Range("G83").Value = "mat2"

$ws.Range("A84").Value = "Szerda"
$ws.Range("B84").Value = 8
$ws.Range("C84").Value = "13.B"
$ws.Range("D84").Value = "Adatbázis-kezelés"
$ws.Range("E84").Value = "HZ"
$ws.Range("F84").Value = 23
$ws.Range("G84").Value = "mat1"

$ws.Range("A85").Value = "Szerda"
$ws.Range("B85").Value = 8
$ws.Range("C85").Value = "13.B"
$ws.Range("D85").Value = "Aalkf."
$ws.Range("E85").Value = "VT"
$ws.Range("F85").Value = 30
$ws.Range("G85").Value = "mat2"

$ws.Range("A86").Value = "Csütörtök"
$ws.Range("B86").Value = 4
$ws.Range("C86").Value = "13.B"
$ws.Range("D86").Value = "Aalkf."
$ws.Range("E86").Value = "VT"
$ws.Range("F86").Value = 30
$ws.Range("G86").Value = "mat1"

$ws.Range("A87").Value = "Csütörtök"
$ws.Range("B87").Value = 4
$ws.Range("C87").Value = "13.B"
$ws.Range("D87").Value = "Frontend"
$ws.Range("E87").Value = "HP"
$ws.Range("F87").Value = 22
$ws.Range("G87").Value = "mat2"

$ws.Range("A88").Value = "Csütörtök"
$ws.Range("B88").Value = 5
$ws.Range("C88").Value = "13.B"
$ws.Range("D88").Value = "Aalkf."
$ws.Range("E88").Value = "VT"
$ws.Range("F88").Value = 30
$ws.Range("G88").Value = "mat1"

$ws.Range("A89").Value = "Csütörtök"
$ws.Range("B89").Value = 5
$ws.Range("C89").Value = "13.B"
$ws.Range("D89").Value = "Frontend"
$ws.Range("E89").Value = "HP"
$ws.Range("F89").Value = 22
$ws.Range("G89").Value = "mat2"

$ws.Range("A90").Value = "Csütörtök"
$ws.Range("B90").Value = 6
$ws.Range("C90").Value = "13.B"
$ws.Range("D90").Value = "Aalkf."
$ws.Range("E90").Value = "VT"
$ws.Range("F90").Value = 30
$ws.Range("G90").Value = "mat1"

$ws.Range("A91").Value = "Csütörtök"
$ws.Range("B91").Value = 6
$ws.Range("C91").Value = "13.B"
$ws.Range("D91").Value = "Frontend"
$ws.Range("E91").Value = "HP"
$ws.Range("F91").Value = 22
$ws.Range("G91").Value = "mat2"

$ws.Range("A92").Value = "Csütörtök"
$ws.Range("B92").Value = 7
$ws.Range("C92").Value = "13.B"
$ws.Range("D92").Value = "Adatbázis-kezelés"
$ws.Range("E92").Value = "HZ"
$ws.Range("F92").Value = 23
$ws.Range("G92").Value = "mat2"

$ws.Range("A93").Value = "Csütörtök"
$ws.Range("B93").Value = 8
$ws.Range("C93").Value = "13.B"
$ws.Range("D93").Value = "Adatbázis-kezelés"
$ws.Range("E93").Value = "HZ"
$ws.Range("F93").Value = 23
$ws.Range("G93").Value = "mat2"

$ws.Range("A94").Value = "Péntek"
$ws.Range("B94").Value = 1
$ws.Range("C94").Value = "13.B"
$ws.Range("D94").Value = "Frontend"
$ws.Range("E94").Value = "ÁA"
$ws.Range("F94").Value = 23
$ws.Range("G94").Value = "mat1"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 90
$ws.Range("J100").Select()

Write-Host "Done"